# Updates the cryptos worksheet with the latest pricing/volume snapshot
# (mirrors "Updated cryptos list ... with GitHub Actions" commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.959.58'
$ws.Range("E2").Value = '  +0.78%  '
$ws.Range("D3").Value = '1.812.71'
$ws.Range("E3").Value = '  +2.13%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.37%  '
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4293'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3695'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.77%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07234'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8641'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.06%  '
$ws.Range("B11").Value = 'WrappedEther'
$ws.Range("C11").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D11").Value = '2.031.84'
$ws.Range("E11").Value = '  +14.31%  '
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.21'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.633'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.22%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.388'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.62%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06893'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.95%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '80.75'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008922'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.68%  '
$ws.Range("E19").Value = '  -0.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.26'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.91%  '
$ws.Range("D21").Value = '26.965.51'
$ws.Range("E21").Value = '  +0.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.195'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.55%  '
$ws.Range("D24").Value = '2.265.20'
$ws.Range("E24").Value = '  +13.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.73'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.880'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.31'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.80%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.226'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.60%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.907'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +16.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '114.98'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08935'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.68%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7446'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.162'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.429'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.798'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.24%  '
$ws.Range("E36").Value = '  +0.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.123'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05213'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01921'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.81%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5082'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.19%  '
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1648'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.40%  '
$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.731'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.426'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.294'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '106.98'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.36'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.004'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4585'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.652'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.80%  '
$ws.Range("E50").Value = '  +0.74%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.809'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.65%  '
